# Update "想去人数" (people interested) counts on both the "展览" sheet
# and the combined "全部类型" sheet, per gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 6450
    $ws.Range("F5").Value = 1015
    $ws.Range("F6").Value = 115
}
